# Update for release to deploy 0.1.1
# - bump Version metadata value 0.1.0 -> 0.1.1
# - bump Date metadata value to the new publish timestamp
# - add a new "Jurisdiction" metadata row (with an empty value) right after
#   "Contact" and before "Description", pushing the remaining metadata rows
#   down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Version: Metadata!B3
$ws.Range("B3").Value = "0.1.1"

# 2) Date: Metadata!B8
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# 3) Insert a new row above the current "Description" row (row 11) and give
#    it the same look (border/alignment) as the rest of the data rows by
#    copying the formatting down from the row that follows the insertion.
$ws.Rows.Item(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
